$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.194993615150452
$ws.Range("B1").Value = 2.593064785003662
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.18999719619751
$ws.Range("E1").Value = 1.178826689720154
